$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 95
$ws.Range("I4").Value = 95
$ws.Range("K4").Value = 95
$ws.Range("M4").Value = 19
$ws.Range("H6").Value = 500008.2
$ws.Range("I6").Value = 500008.2
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1500024.6
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -1499912.6
$ws.Range("N6").ClearContents()
$ws.Range("H11").Value = 47.5
$ws.Range("I11").Value = 47.5
$ws.Range("K11").Value = 47.5
$ws.Range("M11").Value = 92.5
$ws.Range("H28").Value = 1274.25
$ws.Range("I28").Value = 1274.25
$ws.Range("K28").Value = 1274.25
$ws.Range("M28").Value = -789.25
$ws.Range("H40").Value = 1800
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 1950
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 1950
$ws.Range("M40").Value = -1325
$ws.Range("N40").Value = -2300
$ws.Range("H43").Value = 1502
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H74").Value = 2991
$ws.Range("J74").Value = 2991
$ws.Range("L74").Value = 2991
$ws.Range("N74").Value = -4863
$ws.Range("H77").Value = 2991
$ws.Range("J77").Value = 2991
$ws.Range("L77").Value = 14955
$ws.Range("N77").Value = -24315
$ws.Range("H86").Value = 5578.8
$ws.Range("J86").Value = 4473.75
$ws.Range("L86").Value = 4473.75
$ws.Range("N86").Value = -6719.75
$ws.Range("H89").Value = 5578.8
$ws.Range("J89").Value = 4473.75
$ws.Range("L89").Value = 22368.75
$ws.Range("N89").Value = -33600.75
$ws.Range("H113").Value = 3601.6667
$ws.Range("I113").Value = 3452.5
$ws.Range("J113").Value = 3900
$ws.Range("K113").Value = 3452.5
$ws.Range("L113").Value = 3900
$ws.Range("M113").Value = -198.5
$ws.Range("N113").Value = -10408

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6267.5884
$ws.Range("I63").Value = 6136.143
$ws.Range("K63").Value = 6136.143
$ws.Range("M63").Value = -5450.143
$ws.Range("H66").Value = 6267.5884
$ws.Range("I66").Value = 6136.143
$ws.Range("K66").Value = 30680.715
$ws.Range("M66").Value = -27248.715

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 24748.25
$ws.Range("J76").Value = 22497.5
$ws.Range("L76").Value = 22497.5
$ws.Range("N76").Value = -23127.5
$ws.Range("H79").Value = 24748.25
$ws.Range("J79").Value = 22497.5
$ws.Range("L79").Value = 22497.5
$ws.Range("N79").Value = -24681.5

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 119.5
$ws.Range("I10").Value = 119.5
$ws.Range("K10").Value = 119.5
$ws.Range("M10").Value = 19.5
$ws.Range("H22").Value = 216.63637
$ws.Range("I22").Value = 219.3
$ws.Range("K22").Value = 219.3
$ws.Range("M22").Value = 130.7
$ws.Range("H31").Value = 1000
$ws.Range("I31").Value = 1000
$ws.Range("K31").Value = 1000
$ws.Range("M31").Value = -705
$ws.Range("H34").Value = 1000
$ws.Range("I34").Value = 1000
$ws.Range("K34").Value = 1000
$ws.Range("M34").Value = -798
$ws.Range("H58").Value = 2085.5454
$ws.Range("I58").Value = 2054.1
$ws.Range("K58").Value = 2054.1
$ws.Range("M58").Value = -1851.1
$ws.Range("H62").Value = 3016
$ws.Range("I62").Value = 2860
$ws.Range("K62").Value = 2860
$ws.Range("M62").Value = -2236
$ws.Range("H65").Value = 3016
$ws.Range("I65").Value = 2860
$ws.Range("K65").Value = 14300
$ws.Range("M65").Value = -11180
$ws.Range("H136").Value = 2085.5454
$ws.Range("I136").Value = 2054.1
$ws.Range("K136").Value = 6162.299999999999
$ws.Range("M136").Value = -3612.299999999999

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 519642.62
$ws.Range("I9").Value = 500
$ws.Range("J9").Value = 545599.75
$ws.Range("K9").Value = 1500
$ws.Range("L9").Value = 1636799.25
$ws.Range("M9").Value = -1276
$ws.Range("N9").Value = -1637247.25
$ws.Range("H13").Value = 1400.3334
$ws.Range("I13").Value = 1400.3334
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 4201.0002
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -4033.0002
$ws.Range("N13").ClearContents()

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3000000
$ws.Range("I11").Value = 3000000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 3000000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -2999861
$ws.Range("N11").ClearContents()
$ws.Range("H98").Value = 29850.5
$ws.Range("J98").Value = 29850.5
$ws.Range("L98").Value = 29850.5
$ws.Range("N98").Value = -35840.5
$ws.Range("H99").Value = 27999
$ws.Range("I99").Value = 27999
$ws.Range("K99").Value = 27999
$ws.Range("M99").Value = -25753
$ws.Range("H136").Value = 21966.334
$ws.Range("J136").Value = 21966.334
$ws.Range("L136").Value = 65899.00199999999
$ws.Range("N136").Value = -70999.00199999999

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4559.8
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -205
$ws.Range("H27").Value = 4559.8
$ws.Range("I27").Value = 500
$ws.Range("K27").Value = 500
$ws.Range("M27").Value = -393
$ws.Range("H55").Value = 1191.3334
$ws.Range("I55").Value = 866
$ws.Range("K55").Value = 866
$ws.Range("M55").Value = -693
$ws.Range("H68").Value = 2782
$ws.Range("I68").Value = 2660.2
$ws.Range("K68").Value = 2660.2
$ws.Range("M68").Value = -1911.2
$ws.Range("H71").Value = 2782
$ws.Range("I71").Value = 2660.2
$ws.Range("K71").Value = 13301
$ws.Range("M71").Value = -9557
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").ClearContents()
$ws.Range("H132").Value = 12250
$ws.Range("I132").Value = 12250
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 36750
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -34220
$ws.Range("N132").ClearContents()

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3333
$ws.Range("I81").Value = 999.5
$ws.Range("J81").Value = 8000
$ws.Range("K81").Value = 1999
$ws.Range("L81").Value = 16000
$ws.Range("M81").Value = -938
$ws.Range("N81").Value = -18122
$ws.Range("H84").Value = 3333
$ws.Range("I84").Value = 999.5
$ws.Range("J84").Value = 8000
$ws.Range("K84").Value = 9995
$ws.Range("L84").Value = 80000
$ws.Range("M84").Value = -4691
$ws.Range("N84").Value = -90608
$ws.Range("H101").Value = 12598.75
$ws.Range("J101").Value = 12598.75
$ws.Range("L101").Value = 12598.75
$ws.Range("N101").Value = -19088.75
$ws.Range("H104").Value = 10666.333
$ws.Range("J104").Value = 10666.333
$ws.Range("L104").Value = 10666.333
$ws.Range("N104").Value = -17654.333
$ws.Range("H122").Value = 1778.8
$ws.Range("I122").Value = 1298
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 3894
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -1444
$ws.Range("N122").Value = -12400
